# Reproduce the author's edit:
#  - Rename the "deuteron" label (column I / "target", rows 2-10) to "d".
#  - Make the header row (A1:K1) bold + centered (new style).
#  - Move the active-cell selection from L7 to H16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "deuteron" -> "d" for every data row (2-10), column I ---
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq "deuteron") {
        $cell.Value = "d"
    }
}

# --- Bold + center the header row ---
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# --- Update the saved selection/active cell ---
$ws.Range("H16").Select() | Out-Null
